$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.184.42"
$ws.Range("E2").Value = "  -4.63%  "

$ws.Range("D3").Value = "'3.294.92"
$ws.Range("E3").Value = "  -5.37%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'564.79"
$ws.Range("E5").Value = "  -3.65%  "

$ws.Range("D6").Value = "'127.19"
$ws.Range("E6").Value = "  -3.79%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'3.294.77"
$ws.Range("E8").Value = "  -5.36%  "

$ws.Range("D9").Value = "'0.476"
$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "  -4.88%  "

$ws.Range("D11").Value = "'0.117"
$ws.Range("E11").Value = "  -4.55%  "

$ws.Range("D12").Value = "'0.372"
$ws.Range("E12").Value = "  -3.84%  "

$ws.Range("D13").Value = "'3.864.71"
$ws.Range("E13").Value = "  -5.15%  "

$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "'3.303.40"
$ws.Range("E15").Value = "  -5.08%  "

$ws.Range("E16").Value = "  -6.35%  "

$ws.Range("D17").Value = "'61.308.87"
$ws.Range("E17").Value = "  -4.38%  "

$ws.Range("D18").Value = "'24.09"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.60"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").Value = "'8.90"
$ws.Range("E21").Value = "  -11.00%  "

$ws.Range("D22").Value = "'352.65"
$ws.Range("E22").Value = "  -8.47%  "

$ws.Range("D23").Value = "'0.550"
$ws.Range("E23").Value = "  -4.60%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'3.430.81"
$ws.Range("E25").Value = "  -5.28%  "

$ws.Range("D26").Value = "'69.00"
$ws.Range("E26").Value = "  -7.66%  "

$ws.Range("E27").Value = "  -5.77%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'7.10"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "'2.10"
$ws.Range("E31").Value = "  -5.95%  "

$ws.Range("D32").Value = "'7.74"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'3.327.93"
$ws.Range("E35").Value = "  -5.23%  "

$ws.Range("D36").Value = "'22.52"
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").Value = "'5.21"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "'6.74"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D39").Value = "'163.21"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").Value = "'1.46"
$ws.Range("E40").Value = "  -3.56%  "

$ws.Range("D41").Value = "'0.0749"
$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").Value = "'4.36"
$ws.Range("E43").Value = "  +0.67%  "

$ws.Range("D44").Value = "'41.07"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").Value = "'0.740"
$ws.Range("E45").Value = "  -7.70%  "

$ws.Range("E46").Value = "  -2.22%  "

$ws.Range("E47").Value = "  -5.16%  "

$ws.Range("D48").Value = "'22.15"
$ws.Range("E48").Value = "  -7.86%  "

$ws.Range("D49").Value = "'6.64"
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("D50").Value = "'0.847"
$ws.Range("E50").Value = "  -8.22%  "

$ws.Range("D51").Value = "'21.02"
$ws.Range("E51").Value = "  +2.03%  "
